$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 6, 7, 8)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "[-, 'MCT-3A-Eletropneumática', -, -]"
    $ws.Range("E$r").Value = "-"
}
